# Update the "Förändrad" (Changed) date in column C for all data rows
# (rows 2 through 181) from 2023-09-06 to 2023-09-08.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newDate = (Get-Date -Year 2023 -Month 9 -Day 8).Date

$ws.Range("C2:C181").Value = $newDate
